$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.422.59'
$ws.Range('E2').Value = '  +4.21%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.502.61'
$ws.Range('E3').Value = '  +4.02%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '586.24'
$ws.Range('E5').Value = '  +2.95%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.72'
$ws.Range('E6').Value = '  +6.40%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  +1.40%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '7.72'
$ws.Range('E9').Value = '  +0.77%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.126'
$ws.Range('E10').Value = '  +4.51%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.400'
$ws.Range('E11').Value = '  +4.83%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.106.05'
$ws.Range('E12').Value = '  +4.14%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '29.87'
$ws.Range('E13').Value = '  +7.62%  '

$ws.Range('E14').Value = '  -0.50%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.518.51'
$ws.Range('E15').Value = '  +4.49%  '

$ws.Range('E16').Value = '  +4.36%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.483.70'
$ws.Range('E17').Value = '  +4.18%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.27'
$ws.Range('E18').Value = '  +3.28%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '14.29'
$ws.Range('E19').Value = '  +5.54%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '9.50'
$ws.Range('E20').Value = '  +7.13%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '394.43'
$ws.Range('E21').Value = '  +3.62%  '

$ws.Range('E22').Value = '  +3.27%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '75.44'
$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('E25').Value = '  +9.26%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.647.68'
$ws.Range('E26').Value = '  +4.14%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.187'
$ws.Range('E27').Value = '  -0.63%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.81'
$ws.Range('E28').Value = '  +9.37%  '

$ws.Range('E29').Value = '  +0.18%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.27'
$ws.Range('E30').Value = '  +5.91%  '

$ws.Range('E31').Value = '  +2.80%  '

$ws.Range('E32').Value = '  +6.51%  '

$ws.Range('E33').Value = '  +0.05%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.87'
$ws.Range('E34').Value = '  +4.00%  '

$ws.Range('B35').Value = 'EnergySwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '32.70'
$ws.Range('E35').Value = '  +29.23%  '

$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.20'
$ws.Range('E36').Value = '  +5.18%  '

$ws.Range('E37').Value = '  +8.99%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '171.99'
$ws.Range('E38').Value = '  +3.27%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.57'
$ws.Range('E39').Value = '  +9.48%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.542.96'
$ws.Range('E40').Value = '  +4.09%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0772'
$ws.Range('E41').Value = '  +1.44%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.804'
$ws.Range('E42').Value = '  +4.04%  '

$ws.Range('E43').Value = '  +7.98%  '

$ws.Range('E44').Value = '  +4.54%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '42.58'
$ws.Range('E45').Value = '  +0.35%  '

$ws.Range('E46').Value = '  +10.31%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.615.18'
$ws.Range('E47').Value = '  +6.76%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '23.82'
$ws.Range('E48').Value = '  +7.47%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.27'
$ws.Range('E49').Value = '  +12.73%  '

$ws.Range('E51').Value = '  +5.27%  '
